$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 932.41895
$ws.Range("J17").Value = 938.17145
$ws.Range("L17").Value = 2814.51435
$ws.Range("N17").Value = -3150.51435

$ws.Range("H19").Value = 2528.9644
$ws.Range("J19").Value = 1465.9286
$ws.Range("L19").Value = 1465.9286
$ws.Range("N19").Value = -1815.9286

$ws.Range("H75").Value = 50000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 50000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 50000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -51872

$ws.Range("H78").Value = 50000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 50000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 150000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -159360

$ws.Range("H98").Value = 1253.1875
$ws.Range("I98").Value = 1253.1875
$ws.Range("K98").Value = 1253.1875
$ws.Range("M98").Value = 244.8125

$ws.Range("H107").Value = 350.53845
$ws.Range("I107").Value = 345.36
$ws.Range("K107").Value = 345.36
$ws.Range("M107").Value = 1574.64

$ws.Range("H112").Value = 1042.3914
$ws.Range("I112").Value = 824.25
$ws.Range("J112").Value = 1063.1666
$ws.Range("K112").Value = 2472.75
$ws.Range("L112").Value = 3189.4998
$ws.Range("M112").Value = -1364.75
$ws.Range("N112").Value = -5405.4998

$ws.Range("H122").Value = 1253.1875
$ws.Range("I122").Value = 1253.1875
$ws.Range("K122").Value = 3759.5625
$ws.Range("M122").Value = -1309.5625

$ws.Range("H132").Value = 47774.55
$ws.Range("I132").Value = 54416.316
$ws.Range("K132").Value = 163248.948
$ws.Range("M132").Value = -160718.948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10639670
$ws.Range("I32").Value = 11629234
$ws.Range("K32").Value = 11629234
$ws.Range("M32").Value = -11628947

$ws.Range("H45").Value = 2970
$ws.Range("I45").Value = 3011.6316
$ws.Range("J45").Value = 2772.25
$ws.Range("K45").Value = 3011.6316
$ws.Range("L45").Value = 2772.25
$ws.Range("M45").Value = -2634.6316
$ws.Range("N45").Value = -3526.25

$ws.Range("H53").Value = 25000
$ws.Range("I53").Value = 25000
$ws.Range("K53").Value = 25000
$ws.Range("M53").Value = -24318

$ws.Range("H75").Value = 70884.836
$ws.Range("J75").Value = 70884.836
$ws.Range("L75").Value = 70884.836
$ws.Range("N75").Value = -72632.836

$ws.Range("H78").Value = 70884.836
$ws.Range("J78").Value = 70884.836
$ws.Range("L78").Value = 212654.508
$ws.Range("N78").Value = -221390.508

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1654.3334
$ws.Range("J20").Value = 1852.5714
$ws.Range("L20").Value = 1852.5714
$ws.Range("N20").Value = -2346.5714

$ws.Range("H94").Value = 885
$ws.Range("I94").Value = 733.125
$ws.Range("J94").Value = 1155
$ws.Range("K94").Value = 733.125
$ws.Range("L94").Value = 1155
$ws.Range("M94").Value = -282.125
$ws.Range("N94").Value = -2057

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 35000
$ws.Range("J55").Value = 35000
$ws.Range("L55").Value = 35000
$ws.Range("N55").Value = -35630

$ws.Range("H69").Value = 152876
$ws.Range("I69").Value = 136834.67
$ws.Range("K69").Value = 136834.67
$ws.Range("M69").Value = -136085.67

$ws.Range("H72").Value = 152876
$ws.Range("I72").Value = 136834.67
$ws.Range("K72").Value = 410504.01
$ws.Range("M72").Value = -406760.01

$ws.Range("H94").Value = 1250.3158
$ws.Range("I94").Value = 619.1429000000001
$ws.Range("J94").Value = 1618.5
$ws.Range("K94").Value = 619.1429000000001
$ws.Range("L94").Value = 1618.5
$ws.Range("M94").Value = -168.1429000000001
$ws.Range("N94").Value = -2520.5

$ws.Range("H99").Value = 3157.2727
$ws.Range("J99").Value = 2359.6
$ws.Range("L99").Value = 2359.6
$ws.Range("N99").Value = -5355.6

$ws.Range("H105").Value = 20370.945
$ws.Range("I105").Value = 23971.867
$ws.Range("K105").Value = 23971.867
$ws.Range("M105").Value = -22224.867

$ws.Range("H122").Value = 2311.75
$ws.Range("I122").Value = 1298.909
$ws.Range("K122").Value = 3896.727
$ws.Range("M122").Value = -1446.727

$ws.Range("H126").Value = 3157.2727
$ws.Range("J126").Value = 2359.6
$ws.Range("L126").Value = 7078.799999999999
$ws.Range("N126").Value = -12018.8

$ws.Range("H132").Value = 27789858
$ws.Range("I132").Value = 34498330
$ws.Range("K132").Value = 103494990
$ws.Range("M132").Value = -103492460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 98198
$ws.Range("J37").Value = 98198
$ws.Range("L37").Value = 294594
$ws.Range("N37").Value = -294818

$ws.Range("H131").Value = 7741.706
$ws.Range("J131").Value = 9586.444
$ws.Range("L131").Value = 28759.332
$ws.Range("N131").Value = -38839.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3021.2727
$ws.Range("I102").Value = 2987.1428
$ws.Range("J102").Value = 3081
$ws.Range("K102").Value = 2987.1428
$ws.Range("L102").Value = 3081
$ws.Range("M102").Value = -1365.1428
$ws.Range("N102").Value = -6325

$ws.Range("H122").Value = 110893
$ws.Range("I122").Value = 202785
$ws.Range("K122").Value = 608355
$ws.Range("M122").Value = -605905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3386.9167
$ws.Range("I7").Value = 3264.3
$ws.Range("K7").Value = 3264.3
$ws.Range("M7").Value = -3152.3

$ws.Range("H22").Value = 787.2222
$ws.Range("I22").Value = 387
$ws.Range("J22").Value = 1416.1428
$ws.Range("K22").Value = 387
$ws.Range("L22").Value = 1416.1428
$ws.Range("M22").Value = -92
$ws.Range("N22").Value = -2006.1428

$ws.Range("H27").Value = 787.2222
$ws.Range("I27").Value = 387
$ws.Range("J27").Value = 1416.1428
$ws.Range("K27").Value = 387
$ws.Range("L27").Value = 1416.1428
$ws.Range("M27").Value = -280
$ws.Range("N27").Value = -1630.1428

$ws.Range("H45").Value = 10041
$ws.Range("I45").Value = 10041
$ws.Range("K45").Value = 10041
$ws.Range("M45").Value = -9634

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H100").Value = 7472.3184
$ws.Range("I100").Value = 2334.5293
$ws.Range("K100").Value = 2334.5293
$ws.Range("M100").Value = -1793.5293

$ws.Range("H122").Value = 5344.3
$ws.Range("I122").Value = 4915.3076
$ws.Range("J122").Value = 6141
$ws.Range("K122").Value = 14745.9228
$ws.Range("L122").Value = 18423
$ws.Range("M122").Value = -12295.9228
$ws.Range("N122").Value = -23323

$ws.Range("H126").Value = 3386.9167
$ws.Range("I126").Value = 3264.3
$ws.Range("K126").Value = 9792.900000000001
$ws.Range("M126").Value = -7322.900000000001

$ws.Range("H132").Value = 1090075.1
$ws.Range("I132").Value = 1514961
$ws.Range("J132").Value = 4255.778
$ws.Range("K132").Value = 4544883
$ws.Range("L132").Value = 12767.334
$ws.Range("M132").Value = -4542353
$ws.Range("N132").Value = -17827.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1651.2307
$ws.Range("I113").Value = 508
$ws.Range("J113").Value = 2794.4614
$ws.Range("K113").Value = 1524
$ws.Range("L113").Value = 8383.3842
$ws.Range("M113").Value = 646
$ws.Range("N113").Value = -12723.3842

$ws.Range("H122").Value = 3072.6
$ws.Range("I122").Value = 2769.0667
$ws.Range("K122").Value = 8307.2001
$ws.Range("M122").Value = -5857.2001
